$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2422.32
$ws.Range("J17").Value = 2422.32
$ws.Range("L17").Value = 7266.960000000001
$ws.Range("N17").Value = -7602.960000000001

$ws.Range("H62").Value = 1921.4
$ws.Range("I62").Value = 2023.7778
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 2023.7778
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -1399.7778
$ws.Range("N62").Value = -2248

$ws.Range("H65").Value = 1921.4
$ws.Range("I65").Value = 2023.7778
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 10118.889
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -6998.889000000001
$ws.Range("N65").Value = -11240

$ws.Range("H111").Value = 556.1429000000001
$ws.Range("I111").Value = 276.33334
$ws.Range("J111").Value = 766
$ws.Range("K111").Value = 829.0000200000001
$ws.Range("L111").Value = 2298
$ws.Range("M111").Value = 2237.99998
$ws.Range("N111").Value = -8432

$ws.Range("H132").Value = 861282.75
$ws.Range("I132").Value = 1782.6666
$ws.Range("J132").Value = 8167033.5
$ws.Range("K132").Value = 5347.9998
$ws.Range("L132").Value = 24501100.5
$ws.Range("M132").Value = -2817.9998
$ws.Range("N132").Value = -24506160.5

$ws.Range("H137").Value = 3328291.8
$ws.Range("I137").Value = 8674701
$ws.Range("J137").Value = 1637.4889
$ws.Range("K137").Value = 26024103
$ws.Range("L137").Value = 4912.4667
$ws.Range("M137").Value = -26021553
$ws.Range("N137").Value = -10012.4667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6558742
$ws.Range("I32").Value = 9290092
$ws.Range("J32").Value = 3502.8
$ws.Range("K32").Value = 9290092
$ws.Range("L32").Value = 3502.8
$ws.Range("M32").Value = -9289805
$ws.Range("N32").Value = -4076.8

$ws.Range("H61").Value = 37112496
$ws.Range("I61").Value = 47667940
$ws.Range("J61").Value = 168435.67
$ws.Range("K61").Value = 47667940
$ws.Range("L61").Value = 168435.67
$ws.Range("M61").Value = -47667728
$ws.Range("N61").Value = -168859.67

$ws.Range("H74").Value = 10081182
$ws.Range("I74").Value = 16733981
$ws.Range("K74").Value = 16733981
$ws.Range("M74").Value = -16733107

$ws.Range("H77").Value = 10081182
$ws.Range("I77").Value = 16733981
$ws.Range("K77").Value = 83669905
$ws.Range("M77").Value = -83665537

$ws.Range("H132").Value = 59005.742
$ws.Range("I132").Value = 41534.44
$ws.Range("J132").Value = 102684
$ws.Range("K132").Value = 124603.32
$ws.Range("L132").Value = 308052
$ws.Range("M132").Value = -122073.32
$ws.Range("N132").Value = -313112

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 37112496
$ws.Range("I136").Value = 47667940
$ws.Range("J136").Value = 168435.67
$ws.Range("K136").Value = 143003820
$ws.Range("L136").Value = 505307.01
$ws.Range("M136").Value = -143001270
$ws.Range("N136").Value = -510407.01

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H139").Value = 59715
$ws.Range("J139").Value = 59715
$ws.Range("L139").Value = 59715
$ws.Range("N139").Value = -69995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2714.8667
$ws.Range("I86").Value = 2380
$ws.Range("J86").Value = 2882.3
$ws.Range("K86").Value = 2380
$ws.Range("L86").Value = 2882.3
$ws.Range("M86").Value = -1257
$ws.Range("N86").Value = -5128.3

$ws.Range("H89").Value = 2714.8667
$ws.Range("I89").Value = 2380
$ws.Range("J89").Value = 2882.3
$ws.Range("K89").Value = 11900
$ws.Range("L89").Value = 14411.5
$ws.Range("M89").Value = -6284
$ws.Range("N89").Value = -25643.5

$ws.Range("H134").Value = 5291.1177
$ws.Range("I134").Value = 4482.8
$ws.Range("J134").Value = 7536.4443
$ws.Range("K134").Value = 13448.4
$ws.Range("L134").Value = 22609.3329
$ws.Range("M134").Value = -10913.4
$ws.Range("N134").Value = -27679.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2079
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -10885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 827.0353
$ws.Range("I107").Value = 386.48837
$ws.Range("J107").Value = 1278.0714
$ws.Range("K107").Value = 1159.46511
$ws.Range("L107").Value = 3834.2142
$ws.Range("M107").Value = 760.5348900000001
$ws.Range("N107").Value = -7674.2142

$ws.Range("H129").Value = 2084705.4
$ws.Range("I129").Value = 708.3333
$ws.Range("J129").Value = 6413007
$ws.Range("K129").Value = 2124.9999
$ws.Range("L129").Value = 19239021
$ws.Range("M129").Value = 2875.0001
$ws.Range("N129").Value = -19249021

$ws.Range("H131").Value = 923.9737
$ws.Range("J131").Value = 982.1515000000001
$ws.Range("L131").Value = 2946.4545
$ws.Range("N131").Value = -13026.4545

$ws.Range("H132").Value = 2486.6
$ws.Range("I132").Value = 2121.5
$ws.Range("J132").Value = 2730
$ws.Range("K132").Value = 19093.5
$ws.Range("L132").Value = 24570
$ws.Range("M132").Value = -16563.5
$ws.Range("N132").Value = -29630

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1534
$ws.Range("I126").Value = 1500.8572
$ws.Range("K126").Value = 4502.571599999999
$ws.Range("M126").Value = -2032.571599999999

$ws.Range("H132").Value = 47768.297
$ws.Range("I132").Value = 45807.13
$ws.Range("J132").Value = 49916.24
$ws.Range("K132").Value = 137421.39
$ws.Range("L132").Value = 149748.72
$ws.Range("M132").Value = -134891.39
$ws.Range("N132").Value = -154808.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1190.6111
$ws.Range("I107").Value = 1252
$ws.Range("J107").Value = 1094.1428
$ws.Range("K107").Value = 3756
$ws.Range("L107").Value = 3282.4284
$ws.Range("M107").Value = -1836
$ws.Range("N107").Value = -7122.428400000001

$ws.Range("H122").Value = 2573.3044
$ws.Range("I122").Value = 2062.9092
$ws.Range("J122").Value = 3041.1667
$ws.Range("K122").Value = 6188.7276
$ws.Range("L122").Value = 9123.500100000001
$ws.Range("M122").Value = -3738.7276
$ws.Range("N122").Value = -14023.5001

$ws.Range("H132").Value = 120605.94
$ws.Range("I132").Value = 102530.4
$ws.Range("J132").Value = 146428.14
$ws.Range("K132").Value = 307591.2
$ws.Range("L132").Value = 439284.42
$ws.Range("M132").Value = -305061.2
$ws.Range("N132").Value = -444344.42
